$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: SIM replicate 1 - buffer
$ws.Range("D7").Value = "1 Injection SIM.inj"

# Row 8: SIM replicate 1 - receptor
$ws.Range("C8").Value = "Plates Quick.setup"
$ws.Range("D8").Value = "1 Injection SIM.inj"

# Row 9: HDR replicate 1 - buffer
$ws.Range("D9").Value = "ChoderaHDR.inj"

# Row 10: HDR replicate 1 - receptor
$ws.Range("C10").Value = "Plates Quick.setup"
$ws.Range("D10").Value = "ChoderaHDR.inj"

# Row 13: SIM replicate 2(ish) - buffer
$ws.Range("D13").Value = "1 Injection SIM.inj"

# Row 14: SIM replicate 2(ish) - receptor
$ws.Range("C14").Value = "Plates Quick.setup"
$ws.Range("D14").Value = "1 Injection SIM.inj"

# Row 15: HDR replicate 2(ish) - buffer
$ws.Range("D15").Value = "ChoderaHDR.inj"

# Row 16: HDR replicate 2(ish) - receptor
$ws.Range("C16").Value = "Plates Quick.setup"
$ws.Range("D16").Value = "ChoderaHDR.inj"
